# Rename sheets: "TP" -> "LL" (lesion-level / LL ratings), "FP" -> "NL" (non-lesion / NL ratings)
# and update their header-row rating labels accordingly; then drop the
# now-unused Paradigm/FROC/FCTRL columns (D:F) from the TRUTH sheet and
# leave the selection/active-sheet state the way the author left it.

$wb = $excel.ActiveWorkbook

$wsLL = $wb.Worksheets.Item("TP")
$wsLL.Name = "LL"
$wsLL.Range("E1").Value = "LLRating"

$wsNL = $wb.Worksheets.Item("FP")
$wsNL.Name = "NL"
$wsNL.Range("D1").Value = "NLRating"

$wsTruth = $wb.Worksheets.Item("TRUTH")
$wsTruth.Columns("D:F").Delete()

# Restore the per-sheet selections as left by the author.
$wsLL.Range("E1").Select()
$wsTruth.Range("D1:F1048576").Select()

# NL ends up the active/visible sheet when the workbook was last saved.
$wsNL.Activate()
$wsNL.Range("G19").Select()
